$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Re-style the three tables (slides 14-16) from the plain custom table
#    style to the built-in table-style GUID.
# ---------------------------------------------------------------------------
$newTableStyleId = "{52445074-9707-42B3-A11E-902BADFE7834}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId, $true)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Swap the theme colour scheme that is actually in effect for the deck
#    (the "Integral" / Red-Violet palette) for the plain default "Office"
#    palette - i.e. the colours the Office Theme used elsewhere in the
#    package already carries.
#    RGB() packs bytes as 0x00BBGGRR, matching the OLE colour values
#    returned/consumed by ColorScheme.Colors(n).RGB.
# ---------------------------------------------------------------------------
$officeColors = @{
    1  = 0x000000   # dk1      - 000000
    2  = 0xFFFFFF   # lt1      - FFFFFF
    3  = 0x6A5444   # dk2      - 44546A
    4  = 0xE6E6E7   # lt2      - E7E6E6
    5  = 0xD59B5B   # accent1  - 5B9BD5
    6  = 0x317DED   # accent2  - ED7D31
    7  = 0xA5A5A5   # accent3  - A5A5A5
    8  = 0x00C0FF   # accent4  - FFC000
    9  = 0xC47244   # accent5  - 4472C4
    10 = 0x47AD70   # accent6  - 70AD47
    11 = 0xC16305   # hlink    - 0563C1
    12 = 0x724F95   # folHlink - 954F72
}

$colorScheme = $p.SlideMaster.ColorScheme
foreach ($idx in $officeColors.Keys) {
    $colorScheme.Colors($idx).RGB = $officeColors[$idx]
}
